# Auto-generated Excel COM-interop script to apply market-data value updates
# to the Asura_Profits workbook (per-sheet Leve profit tables).
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value2 = 1170.8889
$ws.Range("I28").Value2 = 76
$ws.Range("J28").Value2 = 5003
$ws.Range("K28").Value2 = 76
$ws.Range("L28").Value2 = 5003
$ws.Range("M28").Value2 = 409
$ws.Range("N28").Value2 = -5973
$ws.Range("H31").Value2 = 21597.545
$ws.Range("I31").Value2 = 23557.3
$ws.Range("K31").Value2 = 70671.89999999999
$ws.Range("M31").Value2 = -70441.89999999999
$ws.Range("H53").Value2 = 186.70589
$ws.Range("J53").Value2 = 243.3
$ws.Range("L53").Value2 = 243.3
$ws.Range("N53").Value2 = -1517.3
$ws.Range("H55").Value2 = 412.85715
$ws.Range("I55").Value2 = 101
$ws.Range("J55").Value2 = 464.83334
$ws.Range("K55").Value2 = 101
$ws.Range("L55").Value2 = 464.83334
$ws.Range("M55").Value2 = 113
$ws.Range("N55").Value2 = -892.83334
$ws.Range("H103").Value2 = 700
$ws.Range("I103").Value2 = 0
$ws.Range("J103").Value2 = 700
$ws.Range("K103").Value2 = 0
$ws.Range("L103").Value2 = 2100
$ws.Range("M103").ClearContents() | Out-Null
$ws.Range("N103").Value2 = -3272
$ws.Range("H129").Value2 = 1187.2449
$ws.Range("I129").Value2 = 566.75
$ws.Range("J129").Value2 = 1242.4
$ws.Range("K129").Value2 = 1700.25
$ws.Range("L129").Value2 = 3727.2
$ws.Range("M129").Value2 = 3299.75
$ws.Range("N129").Value2 = -13727.2
$ws.Range("H138").Value2 = 1807.2979
$ws.Range("I138").Value2 = 1380.9032
$ws.Range("J138").Value2 = 2633.4375
$ws.Range("K138").Value2 = 4142.7096
$ws.Range("L138").Value2 = 7900.3125
$ws.Range("M138").Value2 = 997.2903999999999
$ws.Range("N138").Value2 = -18180.3125

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value2 = 17213.842
$ws.Range("I32").Value2 = 22181.656
$ws.Range("K32").Value2 = 22181.656
$ws.Range("M32").Value2 = -21894.656
$ws.Range("H122").Value2 = 5214.9
$ws.Range("I122").Value2 = 6364.8887
$ws.Range("K122").Value2 = 19094.6661
$ws.Range("M122").Value2 = -16644.6661
$ws.Range("H123").Value2 = 24125.4
$ws.Range("J123").Value2 = 24125.4
$ws.Range("L123").Value2 = 24125.4
$ws.Range("N123").Value2 = -33925.4
$ws.Range("H132").Value2 = 1730.2
$ws.Range("I132").Value2 = 1139.2307
$ws.Range("K132").Value2 = 3417.6921
$ws.Range("M132").Value2 = -887.6921000000002

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value2 = 127515.5
$ws.Range("I86").Value2 = 3541.2
$ws.Range("J86").Value2 = 334139.34
$ws.Range("K86").Value2 = 3541.2
$ws.Range("L86").Value2 = 334139.34
$ws.Range("M86").Value2 = -2418.2
$ws.Range("N86").Value2 = -336385.34
$ws.Range("H89").Value2 = 127515.5
$ws.Range("I89").Value2 = 3541.2
$ws.Range("J89").Value2 = 334139.34
$ws.Range("K89").Value2 = 17706
$ws.Range("L89").Value2 = 1670696.7
$ws.Range("M89").Value2 = -12090
$ws.Range("N89").Value2 = -1681928.7

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H23").Value2 = 43300
$ws.Range("J23").Value2 = 0
$ws.Range("L23").Value2 = 0
$ws.Range("N23").ClearContents() | Out-Null
$ws.Range("H27").Value2 = 43300
$ws.Range("J27").Value2 = 0
$ws.Range("L27").Value2 = 0
$ws.Range("N27").ClearContents() | Out-Null
$ws.Range("H31").Value2 = 1558.8837
$ws.Range("I31").Value2 = 1271.1892
$ws.Range("J31").Value2 = 3333
$ws.Range("K31").Value2 = 1271.1892
$ws.Range("L31").Value2 = 3333
$ws.Range("M31").Value2 = -976.1892
$ws.Range("N31").Value2 = -3923
$ws.Range("H34").Value2 = 1558.8837
$ws.Range("I34").Value2 = 1271.1892
$ws.Range("J34").Value2 = 3333
$ws.Range("K34").Value2 = 1271.1892
$ws.Range("L34").Value2 = 3333
$ws.Range("M34").Value2 = -1069.1892
$ws.Range("N34").Value2 = -3737
$ws.Range("H131").Value2 = 24725
$ws.Range("J131").Value2 = 24725
$ws.Range("L131").Value2 = 24725
$ws.Range("N131").Value2 = -34805
$ws.Range("H134").Value2 = 2401.8484
$ws.Range("I134").Value2 = 1752.5
$ws.Range("J134").Value2 = 3400.8462
$ws.Range("K134").Value2 = 5257.5
$ws.Range("L134").Value2 = 10202.5386
$ws.Range("M134").Value2 = -2722.5
$ws.Range("N134").Value2 = -15272.5386

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value2 = 649.8125
$ws.Range("I113").Value2 = 599.65
$ws.Range("J113").Value2 = 733.4167
$ws.Range("K113").Value2 = 1798.95
$ws.Range("L113").Value2 = 2200.2501
$ws.Range("M113").Value2 = 371.0500000000002
$ws.Range("N113").Value2 = -6540.2501
$ws.Range("H122").Value2 = 205.25
$ws.Range("I122").Value2 = 235
$ws.Range("J122").Value2 = 175.5
$ws.Range("K122").Value2 = 2115
$ws.Range("L122").Value2 = 1579.5
$ws.Range("M122").Value2 = 335
$ws.Range("N122").Value2 = -6479.5
$ws.Range("H134").Value2 = 3455.9355
$ws.Range("I134").Value2 = 1791.6842
$ws.Range("J134").Value2 = 6091
$ws.Range("K134").Value2 = 5375.0526
$ws.Range("L134").Value2 = 18273
$ws.Range("M134").Value2 = -305.0526
$ws.Range("N134").Value2 = -28413

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value2 = 50000
$ws.Range("J18").Value2 = 0
$ws.Range("L18").Value2 = 0
$ws.Range("N18").ClearContents() | Out-Null
$ws.Range("H29").Value2 = 175000
$ws.Range("I29").Value2 = 175000
$ws.Range("J29").Value2 = 0
$ws.Range("K29").Value2 = 175000
$ws.Range("L29").Value2 = 0
$ws.Range("M29").Value2 = -174710
$ws.Range("N29").ClearContents() | Out-Null
$ws.Range("H52").Value2 = 30000
$ws.Range("J52").Value2 = 30000
$ws.Range("L52").Value2 = 30000
$ws.Range("N52").Value2 = -30518
$ws.Range("H70").Value2 = 5456.28
$ws.Range("I70").Value2 = 5412.9375
$ws.Range("J70").Value2 = 5533.3335
$ws.Range("K70").Value2 = 5412.9375
$ws.Range("L70").Value2 = 5533.3335
$ws.Range("M70").Value2 = -5142.9375
$ws.Range("N70").Value2 = -6073.3335
$ws.Range("H73").Value2 = 5456.28
$ws.Range("I73").Value2 = 5412.9375
$ws.Range("J73").Value2 = 5533.3335
$ws.Range("K73").Value2 = 5412.9375
$ws.Range("L73").Value2 = 5533.3335
$ws.Range("M73").Value2 = -4476.9375
$ws.Range("N73").Value2 = -7405.3335
$ws.Range("H122").Value2 = 2935.0386
$ws.Range("I122").Value2 = 2099.7693
$ws.Range("J122").Value2 = 3770.3076
$ws.Range("K122").Value2 = 6299.3079
$ws.Range("L122").Value2 = 11310.9228
$ws.Range("M122").Value2 = -3849.3079
$ws.Range("N122").Value2 = -16210.9228
$ws.Range("H123").Value2 = 8923.272000000001
$ws.Range("J123").Value2 = 8923.272000000001
$ws.Range("L123").Value2 = 8923.272000000001
$ws.Range("N123").Value2 = -13823.272
$ws.Range("H131").Value2 = 47658
$ws.Range("J131").Value2 = 47658
$ws.Range("L131").Value2 = 47658
$ws.Range("N131").Value2 = -57738
$ws.Range("H132").Value2 = 2099.238
$ws.Range("I132").Value2 = 1504.2941
$ws.Range("J132").Value2 = 4627.75
$ws.Range("K132").Value2 = 4512.8823
$ws.Range("L132").Value2 = 13883.25
$ws.Range("M132").Value2 = -1982.8823
$ws.Range("N132").Value2 = -18943.25

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value2 = 2280.8
$ws.Range("I61").Value2 = 2801.3333
$ws.Range("J61").Value2 = 1500
$ws.Range("K61").Value2 = 2801.3333
$ws.Range("L61").Value2 = 1500
$ws.Range("M61").Value2 = -2599.3333
$ws.Range("N61").Value2 = -1904
$ws.Range("H113").Value2 = 2280.8
$ws.Range("I113").Value2 = 2801.3333
$ws.Range("J113").Value2 = 1500
$ws.Range("K113").Value2 = 2801.3333
$ws.Range("L113").Value2 = 1500
$ws.Range("M113").Value2 = -631.3332999999998
$ws.Range("N113").Value2 = -5840
$ws.Range("H122").Value2 = 25006000
$ws.Range("I122").Value2 = 6666.6665
$ws.Range("J122").Value2 = 100004000
$ws.Range("K122").Value2 = 19999.9995
$ws.Range("L122").Value2 = 300012000
$ws.Range("M122").Value2 = -17549.9995
$ws.Range("N122").Value2 = -300016900
$ws.Range("H132").Value2 = 4444.7295
$ws.Range("I132").Value2 = 4484
$ws.Range("K132").Value2 = 13452
$ws.Range("M132").Value2 = -10922

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H17").Value2 = 17951.666
$ws.Range("I17").Value2 = 17951.666
$ws.Range("K17").Value2 = 17951.666
$ws.Range("M17").Value2 = -17779.666
$ws.Range("H29").Value2 = 9125
$ws.Range("I29").Value2 = 4000
$ws.Range("J29").Value2 = 10833.333
$ws.Range("K29").Value2 = 4000
$ws.Range("L29").Value2 = 10833.333
$ws.Range("M29").Value2 = -3710
$ws.Range("N29").Value2 = -11413.333
$ws.Range("H51").Value2 = 10500
$ws.Range("I51").Value2 = 6000
$ws.Range("J51").Value2 = 15000
$ws.Range("K51").Value2 = 6000
$ws.Range("L51").Value2 = 15000
$ws.Range("M51").Value2 = -5490
$ws.Range("N51").Value2 = -16020
$ws.Range("H64").Value2 = 35000
$ws.Range("J64").Value2 = 35000
$ws.Range("L64").Value2 = 35000
$ws.Range("N64").Value2 = -35496
$ws.Range("H67").Value2 = 35000
$ws.Range("J67").Value2 = 35000
$ws.Range("L67").Value2 = 35000
$ws.Range("N67").Value2 = -36716
$ws.Range("H86").Value2 = 200081.25
$ws.Range("J86").Value2 = 200081.25
$ws.Range("L86").Value2 = 200081.25
$ws.Range("N86").Value2 = -202327.25
$ws.Range("H89").Value2 = 200081.25
$ws.Range("J89").Value2 = 200081.25
$ws.Range("L89").Value2 = 1000406.25
$ws.Range("N89").Value2 = -1011638.25
$ws.Range("H122").Value2 = 100005064
$ws.Range("I122").Value2 = 200004130
$ws.Range("J122").Value2 = 6001
$ws.Range("K122").Value2 = 600012390
$ws.Range("L122").Value2 = 18003
$ws.Range("M122").Value2 = -600009940
$ws.Range("N122").Value2 = -22903
$ws.Range("H123").Value2 = 37479.81
$ws.Range("J123").Value2 = 37479.81
$ws.Range("L123").Value2 = 37479.81
$ws.Range("N123").Value2 = -47279.81
$ws.Range("H126").Value2 = 12212.429
$ws.Range("I126").Value2 = 13414.5
$ws.Range("J126").Value2 = 5000
$ws.Range("K126").Value2 = 40243.5
$ws.Range("L126").Value2 = 15000
$ws.Range("M126").Value2 = -37773.5
$ws.Range("N126").Value2 = -19940
$ws.Range("H132").Value2 = 2109.2263
$ws.Range("I132").Value2 = 1234.4333
$ws.Range("J132").Value2 = 3250.261
$ws.Range("K132").Value2 = 3703.2999
$ws.Range("L132").Value2 = 9750.782999999999
$ws.Range("M132").Value2 = -1173.2999
$ws.Range("N132").Value2 = -14810.783
$ws.Range("H136").Value2 = 1354.1936
$ws.Range("I136").Value2 = 1225.2222
$ws.Range("J136").Value2 = 2224.75
$ws.Range("K136").Value2 = 3675.6666
$ws.Range("L136").Value2 = 6674.25
$ws.Range("M136").Value2 = -1125.6666
$ws.Range("N136").Value2 = -11774.25
